$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.680.01'
$ws.Range("E2").Value = '  +2.88%  '

# Row 3
$ws.Range("D3").Value = '2.946.63'
$ws.Range("E3").Value = '  +2.11%  '

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = "'591.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '

# Row 6
$ws.Range("D6").Value = "'148.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.43%  '

# Row 7
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '2.946.41'
$ws.Range("E8").Value = '  +2.12%  '

# Row 9
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = "'0.508"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.06%  '

# Row 10
$ws.Range("D10").Value = "'7.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.01%  '

# Row 11
$ws.Range("E11").Value = '  +10.19%  '

# Row 12
$ws.Range("E12").Value = '  +2.58%  '

# Row 13
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.32%  '

# Row 14
$ws.Range("D14").Value = "'32.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.88%  '

# Row 16
$ws.Range("D16").Value = '3.435.43'
$ws.Range("E16").Value = '  +2.10%  '

# Row 17
$ws.Range("D17").Value = '62.623.65'
$ws.Range("E17").Value = '  +2.97%  '

# Row 18
$ws.Range("E18").Value = '  +2.91%  '

# Row 19
$ws.Range("D19").Value = '2.945.62'
$ws.Range("E19").Value = '  +2.82%  '

# Row 20
$ws.Range("D20").Value = "'438.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.82%  '

# Row 21
$ws.Range("D21").Value = "'13.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.40%  '

# Row 22
$ws.Range("D22").Value = "'0.666"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.00%  '

# Row 23
$ws.Range("E23").Value = '  +1.02%  '

# Row 24
$ws.Range("E24").Value = '  +8.78%  '

# Row 25
$ws.Range("D25").Value = "'80.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.74%  '

# Row 26
$ws.Range("D26").Value = "'11.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.36%  '

# Row 27
$ws.Range("D27").Value = "'2.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.62%  '

# Row 28
$ws.Range("E28").Value = '  +0.07%  '

# Row 29
$ws.Range("D29").Value = "'7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.00%  '

# Row 30
$ws.Range("D30").Value = "'0.0000104"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +24.00%  '

# Row 31
$ws.Range("D31").Value = "'2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.72%  '

# Row 32
$ws.Range("E32").Value = '  +5.39%  '

# Row 33
$ws.Range("E33").Value = '  +6.47%  '

# Row 34
$ws.Range("D34").Value = "'26.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.44%  '

# Row 35
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("E36").Value = '  +2.07%  '

# Row 37
$ws.Range("D37").Value = "'3.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.98%  '

# Row 38
$ws.Range("D38").Value = "'5.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.09%  '

# Row 39
$ws.Range("D39").Value = "'49.67"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.49%  '

# Row 40
$ws.Range("D40").Value = "'2.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.11%  '

# Row 41
$ws.Range("E41").Value = '  +1.21%  '

# Row 42
$ws.Range("D42").Value = "'0.117"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.90%  '

# Row 43
$ws.Range("D43").Value = "'0.279"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.12%  '

# Row 44
$ws.Range("D44").Value = "'40.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.18%  '

# Row 45
$ws.Range("D45").Value = '2.707.21'
$ws.Range("E45").Value = '  +1.88%  '

# Row 46
$ws.Range("D46").Value = "'135.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.45%  '

# Row 47
$ws.Range("E47").Value = '  +4.14%  '

# Row 48
$ws.Range("D48").Value = "'357.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.00%  '

# Row 50
$ws.Range("E50").Value = '  +2.49%  '

# Row 51
$ws.Range("D51").Value = "'22.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.56%  '
